$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($val -eq "backup@backdoor.com, System, system") {
        $cell.Value = "system, backup@backdoor.com, System"
    } elseif ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    } elseif ($val -eq "System, admin@admin.com") {
        $cell.Value = "admin@admin.com, System"
    }
}
